$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as an exact text string (preserves formatting such as
# trailing zeros, e.g. "5.182" or "0.0002000") even though the text looks numeric.
# We temporarily force a Text number format so Excel does not coerce the string
# into a floating point number, then restore the cell style to Normal so no
# visible/formatting side effects remain.
function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue "D2" "245.72"
Set-TextValue "D3" "23.89"
Set-TextValue "D4" "5.182"
Set-TextValue "D5" "0.05741"
Set-TextValue "D6" "6.479"
Set-TextValue "D7" "3.151"
Set-TextValue "D8" "0.8132"
Set-TextValue "D9" "0.8552"
Set-TextValue "D10" "0.1377"
Set-TextValue "D11" "0.06956"
Set-TextValue "D12" "0.03190"
Set-TextValue "D13" "0.02889"
Set-TextValue "D14" "0.09376"
Set-TextValue "D15" "3.820"
Set-TextValue "D16" "0.001528"
Set-TextValue "D17" "0.04698"
Set-TextValue "D18" "0.0005998"
Set-TextValue "D19" "0.006186"
Set-TextValue "E19" "18TigerCashTCHBestin24h"
Set-TextValue "D20" "0.001241"
Set-TextValue "D21" "0.004786"
Set-TextValue "E21" "20HotbitTokenHTB"
Set-TextValue "D22" "0.00008494"
Set-TextValue "D24" "2.144"
Set-TextValue "D25" "0.3203"
Set-TextValue "D27" "0.1328"
Set-TextValue "D28" "0.0002331"
Set-TextValue "D40" "0.03699"
Set-TextValue "D42" "0.1055"
Set-TextValue "D43" "0.002205"
Set-TextValue "D44" "0.007800"
Set-TextValue "D45" "0.00005497"
Set-TextValue "D47" "0.3999"
Set-TextValue "D48" "0.002738"
Set-TextValue "D49" "0.00002100"
Set-TextValue "D50" "0.0002000"
